$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(5, 11, 15, 18, 20, 22, 27, 29, 33, 37, 38, 39, 42, 45, 46)

foreach ($r in $rows) {
    $rangeAddr = "H" + $r + ":K" + $r
    $ws.Range($rangeAddr).Value = "Transporte"
}
